# Mise à jour de l'application
# Adds two new training-session date columns (AA = 11/08/2025, AB = 13/08/2025)
# to the attendance sheet and records each player's attendance ("P" = présent,
# "RH" = retour vers un autre club / indisponible) for the two new sessions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: the two new session dates, formatted like the existing ones ---
$ws.Range("AA1").Value = 45880
$ws.Range("AB1").Value = 45882
$ws.Range("Z1").Copy()
$ws.Range("AA1:AB1").PasteSpecial(-4122)

# --- Rows that should be marked "RH" instead of "P" for the two new sessions ---
$rhRows = @(17, 27)

# --- Data rows 2-27: one cell per session, formatted like the rest of the row ---
for ($r = 2; $r -le 27; $r++) {
    $value = "P"
    if ($rhRows -contains $r) {
        $value = "RH"
    }

    $ws.Range("AA$r").Value = $value
    $ws.Range("AB$r").Value = $value

    $ws.Range("Z$r").Copy()
    $ws.Range("AA" + $r + ":AB" + $r).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
[void]$ws.Range("AC23").Select()
